# Refresh the crypto price/volume snapshot (cols D "Price", E "Volume(1h)")
# for rows 2-51, per the scraper run on Sat Nov  2 13:12:00 UTC 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.518.28'
$ws.Range('E2').Value = '  -1.04%  '
$ws.Range('D3').Value = '2.491.36'
$ws.Range('E3').Value = '  -1.69%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '569.47'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.72%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '165.52'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.45%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  -1.47%  '
$ws.Range('D9').Value = '2.489.45'
$ws.Range('E9').Value = '  -1.74%  '
$ws.Range('E10').Value = '  -3.13%  '
$ws.Range('E11').Value = '  -0.41%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.354'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.91'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.20%  '
$ws.Range('D14').Value = '2.942.02'
$ws.Range('E14').Value = '  -1.52%  '
$ws.Range('D15').Value = '69.443.92'
$ws.Range('E15').Value = '  -0.98%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000175'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.34%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '24.29'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -4.09%  '
$ws.Range('D18').Value = '2.487.24'
$ws.Range('E18').Value = '  -1.83%  '
$ws.Range('E19').Value = '  -2.09%  '
$ws.Range('E20').Value = '  -7.07%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '346.91'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.84%  '
$ws.Range('E22').Value = '  -1.92%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.92'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -3.46%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '70.34'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.53%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.88'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -3.64%  '
$ws.Range('D27').Value = '2.615.56'
$ws.Range('E27').Value = '  -2.23%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.64'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -5.17%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.995'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.01%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.80'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.31%  '
$ws.Range('D31').Value = '0.0₃0877'
$ws.Range('E31').Value = '  -4.16%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '452.16'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -3.30%  '
$ws.Range('E33').Value = '  -5.93%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.999'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('E35').Value = '  -3.21%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '155.90'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.55%  '
$ws.Range('E37').Value = '  -4.63%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '19.03'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '18.20'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -2.11%  '
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('E41').Value = '  -2.36%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '4.62'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -4.04%  '
$ws.Range('E43').Value = '  -1.39%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '37.99'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.75%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.16'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -7.63%  '
$ws.Range('E46').Value = '  -8.31%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '139.65'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.47%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.44'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.29%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.513'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -4.09%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0730'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.16%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.574'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.83%  '
